$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the content of rows 54-58 down to rows 55-59 (no structural insert; keep
# existing per-row alternating style formatting intact), then place the new
# "SCES" / "Sem.Césure" entry at row 54 and re-create the trailing blank
# spacer row at 60 (it used to be row 59).
for ($r = 58; $r -ge 54; $r--) {
    $ws.Cells.Item($r + 1, 2).Value2 = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r + 1, 3).Value2 = $ws.Cells.Item($r, 3).Value2
}

$ws.Cells.Item(54, 2).Value2 = "SCES"
$ws.Cells.Item(54, 3).Value2 = "Sem.Césure"

# Row 59 now holds real data (it used to be the blank spacer row), so give it
# the same alternating zebra-stripe cell formatting as the other odd data
# rows (e.g. row 57), and the standard data-row height.
$ws.Range("B57:C57").Copy()
$ws.Range("B59:C59").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Rows.Item(59).RowHeight = 19.7321

# Row 60 becomes the new trailing blank spacer row, at the taller height the
# spacer row used to have.
$ws.Rows.Item(60).RowHeight = 28.7982
